{"js": "// Update the date line and every \"A\u00d7B=\" practice problem in the\n// three-digit-by-one-digit multiplication worksheet.\n//\n// The worksheet's single table has 20 rows \u00d7 5 columns, but only rows\n// 0, 4, 9, 14 and 19 actually contain problem text (the rows in\n// between are blank answer rows). We update each populated cell by\n// position (row/col) rather than by a global text search, because a\n// couple of the new values (\"286\u00d76=\") coincide with old values that\n// live elsewhere in the table \u2014 a plain find/replace pass could hit\n// the wrong (already-updated) cell.\n\nconst body = context.document.body;\n\n// 1) Title / date paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nconst oldTitle = \"2026-01-23 Friday\";\nconst newTitle = \"2026-01-24 Saturday\";\nif (titlePara.text.trim() === oldTitle) {\n  titlePara.getRange().insertText(newTitle, Word.InsertLocation.replace);\n}\n\n// 2) Table cells: (row, col, oldText, newText).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, oldText: \"285\u00d77=\", newText: \"745\u00d75=\" },\n  { row: 0, col: 1, oldText: \"178\u00d77=\", newText: \"837\u00d79=\" },\n  { row: 0, col: 2, oldText: \"246\u00d77=\", newText: \"151\u00d75=\" },\n  { row: 0, col: 3, oldText: \"578\u00d72=\", newText: \"406\u00d72=\" },\n  { row: 0, col: 4, oldText: \"390\u00d73=\", newText: \"206\u00d74=\" },\n  { row: 4, col: 0, oldText: \"700\u00d78=\", newText: \"680\u00d76=\" },\n  { row: 4, col: 1, oldText: \"895\u00d73=\", newText: \"668\u00d73=\" },\n  { row: 4, col: 2, oldText: \"445\u00d79=\", newText: \"179\u00d73=\" },\n  { row: 4, col: 3, oldText: \"667\u00d76=\", newText: \"734\u00d75=\" },\n  { row: 4, col: 4, oldText: \"183\u00d79=\", newText: \"495\u00d75=\" },\n  { row: 9, col: 0, oldText: \"202\u00d77=\", newText: \"155\u00d74=\" },\n  { row: 9, col: 1, oldText: \"102\u00d79=\", newText: \"286\u00d76=\" },\n  { row: 9, col: 2, oldText: \"934\u00d74=\", newText: \"661\u00d76=\" },\n  { row: 9, col: 3, oldText: \"850\u00d79=\", newText: \"832\u00d72=\" },\n  { row: 9, col: 4, oldText: \"936\u00d73=\", newText: \"354\u00d79=\" },\n  { row: 14, col: 0, oldText: \"636\u00d79=\", newText: \"963\u00d73=\" },\n  { row: 14, col: 1, oldText: \"520\u00d74=\", newText: \"905\u00d79=\" },\n  { row: 14, col: 2, oldText: \"972\u00d74=\", newText: \"367\u00d78=\" },\n  { row: 14, col: 3, oldText: \"530\u00d79=\", newText: \"464\u00d78=\" },\n  { row: 14, col: 4, oldText: \"694\u00d77=\", newText: \"631\u00d79=\" },\n  { row: 19, col: 0, oldText: \"121\u00d72=\", newText: \"449\u00d75=\" },\n  { row: 19, col: 1, oldText: \"750\u00d77=\", newText: \"820\u00d78=\" },\n  { row: 19, col: 2, oldText: \"105\u00d77=\", newText: \"345\u00d72=\" },\n  { row: 19, col: 3, oldText: \"228\u00d72=\", newText: \"255\u00d77=\" },\n  { row: 19, col: 4, oldText: \"286\u00d76=\", newText: \"869\u00d76=\" }\n];\n\nconst cells = updates.map((u) => table.getCell(u.row, u.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const cell = cells[i];\n  const expected = updates[i].oldText;\n  if (cell.value.trim() === expected) {\n    cell.body.getRange().insertText(updates[i].newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"A\u00d7B=\" practice problem in the\n# three-digit-by-one-digit multiplication worksheet.\n#\n# The worksheet's single table has 20 rows x 5 columns, but only rows\n# 1, 5, 10, 15 and 20 (1-based) actually contain problem text (the\n# rows in between are blank answer rows). Each populated cell is\n# updated in place by setting Range.Text directly (after trimming the\n# trailing end-of-cell mark with MoveEnd) rather than via a document-\n# wide Find/Replace: a couple of the new values (\"286x6=\") coincide\n# with old values that live elsewhere in the table, so a global\n# find-and-replace pass run after earlier replacements have already\n# landed could re-match and clobber the wrong (already-updated) cell.\n# Scoping the edit to each cell's own Range sidesteps that entirely.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($range, [string]$oldText, [string]$newText) {\n    $range.MoveEnd(1, -1) | Out-Null   # wdCharacter: drop trailing cell/paragraph mark\n    if ($range.Text -eq $oldText) {\n        $range.Text = $newText\n    }\n}\n\n# 1) Title / date paragraph (first paragraph in the body).\n$titleRange = $d.Paragraphs.Item(1).Range\nSet-RangeText $titleRange \"2026-01-23 Friday\" \"2026-01-24 Saturday\"\n\n# 2) Table cells: (row, col, oldText, newText) - 1-based row/col.\n$table = $d.Tables.Item(1)\n\n$updates = @(\n    @(1, 1, \"285\u00d77=\", \"745\u00d75=\"),\n    @(1, 2, \"178\u00d77=\", \"837\u00d79=\"),\n    @(1, 3, \"246\u00d77=\", \"151\u00d75=\"),\n    @(1, 4, \"578\u00d72=\", \"406\u00d72=\"),\n    @(1, 5, \"390\u00d73=\", \"206\u00d74=\"),\n    @(5, 1, \"700\u00d78=\", \"680\u00d76=\"),\n    @(5, 2, \"895\u00d73=\", \"668\u00d73=\"),\n    @(5, 3, \"445\u00d79=\", \"179\u00d73=\"),\n    @(5, 4, \"667\u00d76=\", \"734\u00d75=\"),\n    @(5, 5, \"183\u00d79=\", \"495\u00d75=\"),\n    @(10, 1, \"202\u00d77=\", \"155\u00d74=\"),\n    @(10, 2, \"102\u00d79=\", \"286\u00d76=\"),\n    @(10, 3, \"934\u00d74=\", \"661\u00d76=\"),\n    @(10, 4, \"850\u00d79=\", \"832\u00d72=\"),\n    @(10, 5, \"936\u00d73=\", \"354\u00d79=\"),\n    @(15, 1, \"636\u00d79=\", \"963\u00d73=\"),\n    @(15, 2, \"520\u00d74=\", \"905\u00d79=\"),\n    @(15, 3, \"972\u00d74=\", \"367\u00d78=\"),\n    @(15, 4, \"530\u00d79=\", \"464\u00d78=\"),\n    @(15, 5, \"694\u00d77=\", \"631\u00d79=\"),\n    @(20, 1, \"121\u00d72=\", \"449\u00d75=\"),\n    @(20, 2, \"750\u00d77=\", \"820\u00d78=\"),\n    @(20, 3, \"105\u00d77=\", \"345\u00d72=\"),\n    @(20, 4, \"228\u00d72=\", \"255\u00d77=\"),\n    @(20, 5, \"286\u00d76=\", \"869\u00d76=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $oldText = $u[2]\n    $newText = $u[3]\n    $cellRange = $table.Cell($row, $col).Range\n    Set-RangeText $cellRange $oldText $newText\n}\n"}
